$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aufgb3")
$ws.Activate()

# Clear the AccessCount table contents for the slow MatrixPot data (rows 81-139, columns B:D)
$ws.Range("B81:D139").ClearContents()

# Update the view state to match: scrolled to A62, selection at E87
$ws.Range("E87").Select()
$excel.ActiveWindow.ScrollRow = 62
